# Applies crypto price/volume updates and a few row (coin) reorderings
# as described in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.391.04'
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = '2.649.68'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.38'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.73'
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.57'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.379'
$ws.Range("E11").Value = '  +2.54%  '
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '3.122.52'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.40'
$ws.Range("E14").Value = '  +7.05%  '
$ws.Range("D15").Value = '61.399.76'
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").Value = '2.667.64'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.60'
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.75'
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.18'
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.90'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.524'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.06'
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.58'
$ws.Range("E25").Value = '  +6.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  +2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  +6.42%  '
$ws.Range("D29").Value = '0.0₃0819'
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.82'
$ws.Range("E30").Value = '  +6.87%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.97'
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.03'
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.10'
$ws.Range("E34").Value = '  +12.60%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.67'
$ws.Range("E35").Value = '  +8.80%  '
$ws.Range("E36").Value = '  +8.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.72'
$ws.Range("E37").Value = '  +5.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '341.27'
$ws.Range("E38").Value = '  +9.73%  '
$ws.Range("E39").Value = '  +6.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.907'
$ws.Range("E40").Value = '  +8.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.35'
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.35'
$ws.Range("E42").Value = '  +6.71%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.22'
$ws.Range("E43").Value = '  +5.28%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0577'
$ws.Range("E44").Value = '  +4.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.628'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '20.46'
$ws.Range("E46").Value = '  +3.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.83'
$ws.Range("E47").Value = '  -0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0251'
$ws.Range("E48").Value = '  +4.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0997'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").Value = '2.099.08'
$ws.Range("E51").Value = '  +2.80%  '
